# #5: property aircraft done
# Fix property_category column values that were incorrectly left as "land"
# on the building (建物) and car (汽車) sheets.

$wb = $excel.ActiveWorkbook

# --- 建物 (building) sheet: property_category column is I, rows 2-4 ---
$wsBuilding = $wb.Worksheets.Item("建物")
$wsBuilding.Range("I2").Value = "building"
$wsBuilding.Range("I3").Value = "building"
$wsBuilding.Range("I4").Value = "building"

# --- 汽車 (car) sheet: property_category column is H, row 2 ---
$wsCar = $wb.Worksheets.Item("汽車")
$wsCar.Range("H2").Value = "car"
